# Atualização automática de pedidos - 30/05/2025 08:26
#
# 1) Sheet "Pedidos": mark REQ-006 (row 7) as concluded, stamp last-update
#    info, then append a new order REQ-007 in row 8.
# 2) Sheet "Itens": normalize REQ-006 RACK quantity to a real number, then
#    append the (still empty/untyped) item line for REQ-007 in row 8.

$wb = $excel.ActiveWorkbook

$wsPedidos = $wb.Worksheets.Item("Pedidos")
$wsItens   = $wb.Worksheets.Item("Itens")

# --- Pedidos: update existing REQ-006 row (row 7) ---------------------
$wsPedidos.Range("D7").Value = 1
$wsPedidos.Range("H7").Value = "Concluído"
$wsPedidos.Range("I7").Value = "30/05/2025 07:33"
$wsPedidos.Range("J7").Value = "teste"

# --- Pedidos: append new REQ-007 row (row 8) ---------------------------
$wsPedidos.Range("A8").Value = "REQ-007"
$wsPedidos.Range("B8").Value = "30/05/2025 08:25"
$wsPedidos.Range("C8").Value = "Ford"
$wsPedidos.Range("D8").Value = "'1"
$wsPedidos.Range("E8").Value = "R01-LA-A2"
$wsPedidos.Range("F8").Value = "teste"
$wsPedidos.Range("G8").Value = "'"
$wsPedidos.Range("H8").Value = "Pendente"
$wsPedidos.Range("I8").Value = "'"
$wsPedidos.Range("J8").Value = "'"

# --- Itens: fix REQ-006 seccao to a numeric value (row 7) --------------
$wsItens.Range("D7").Value = 1

# --- Itens: append placeholder line for REQ-007 (row 8) ----------------
$wsItens.Range("A8").Value = "REQ-007"
$wsItens.Range("B8").Value = "'"
$wsItens.Range("C8").Value = "'0"
$wsItens.Range("D8").Value = "'0"
$wsItens.Range("E8").Value = "'0"
$wsItens.Range("F8").Value = 1
